# Apply hybrid bold + color (#2C3E50) highlighting to quantitative metrics
# (percentages, dollar amounts, large numbers) across specific bullet points,
# matching the target OOXML diff.

$d = $word.ActiveDocument

# RGB(0x2C,0x3E,0x50) encoded for Word's BGR COM color integer.
$metricColor = 5258796

function Highlight-Metric {
    param(
        $Paragraph,
        [string]$SearchText
    )
    $r = $Paragraph.Range
    $found = $r.Find.Execute($SearchText, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        $r.Font.Bold = 1
        $r.Font.Color = $metricColor
    }
    return $found
}

# Paragraph 9: "...demographic classification accuracy from 23% to 64%"
$p = $d.Paragraphs(9)
Highlight-Metric $p "23%"
Highlight-Metric $p "64%"

# Paragraph 11: "Achieved 87% prediction accuracy ... industry standard of 71%,
#                reducing polling error margins from ±4.2% to ±2.1%"
$plusMinus = [char]0x00B1
$lowMargin = $plusMinus + "4.2%"
$highMargin = $plusMinus + "2.1%"
$p = $d.Paragraphs(11)
Highlight-Metric $p "87%"
Highlight-Metric $p "71%"
Highlight-Metric $p $lowMargin
Highlight-Metric $p $highMargin

# Paragraph 31: "...analyzed bids from 1,200 vendors..."
$p = $d.Paragraphs(31)
Highlight-Metric $p "1,200"

# Paragraph 46: "...became the $400M Polling Consortium Database ...
#                now valued at $1B+"
$p = $d.Paragraphs(46)
Highlight-Metric $p "`$400M"
Highlight-Metric $p "`$1B"

# Paragraph 63: "Algorithm reduced mapping costs by 73.5%, saving campaigns
#                and organizations $4.7M"
$p = $d.Paragraphs(63)
Highlight-Metric $p "73.5%"
Highlight-Metric $p "`$4.7M"

# Paragraph 65: "Achieved 87% prediction accuracy ... industry standard of 71%"
$p = $d.Paragraphs(65)
Highlight-Metric $p "87%"
Highlight-Metric $p "71%"
